# Apply DASHBOARD.xlsx changes described in the commit diff.

$wb = $excel.ActiveWorkbook
$wsSlides = $wb.Worksheets.Item("Slides")
$wsBatches = $wb.Worksheets.Item("Batches")

# --- Slides sheet ---

# Row 2 (S12 - SCOT-HEART): updated summary text + trailing dash style
$wsSlides.Range("H2").Value = "Figura incluída (KM esquemático) + citação SCOT-HEART 10y corrigida; padding ajustado."
$wsSlides.Range("J2").Value = "-"

# Row 13 (S18 - Imprecisão aplicada ao CLEAR Outcomes)
$wsSlides.Range("F13").Value = "P0"
# Force the date-looking text to stay plain text (column holds strings, not
# real dates) instead of Excel auto-converting "2026-01-23" into a date.
$wsSlides.Range("G13").NumberFormat = "@"
$wsSlides.Range("G13").Value = "2026-01-23"
$wsSlides.Range("G13").ClearFormats()
$wsSlides.Range("H13").Value = "Revisão de consistência visual (mantido layout); segue como base para downgrade por imprecisão"
$wsSlides.Range("J13").Value = "-"

# Row 14 (S19 - Risco de viés RoB 2.0 aplicado ao CLEAR Outcomes)
$wsSlides.Range("F14").Value = "P0"
$wsSlides.Range("G14").NumberFormat = "@"
$wsSlides.Range("G14").Value = "2026-01-23"
$wsSlides.Range("G14").ClearFormats()
$wsSlides.Range("H14").Value = "Decisão final sem emoji; selo ✓ consistente (navy/teal)"
$wsSlides.Range("J14").Value = "-"

# Row 16 (S23 - GRADE: Prognóstico vs Intervenção) collapses to the same
# "DONE" shorthand pattern used by rows 17/19 - only columns A-D stay, E becomes DONE.
$wsSlides.Range("E16").Value = "DONE"
$wsSlides.Range("F16:J16").ClearContents()

# Row 18 (S25 - GRADE: Discriminação do PREVENT)
$wsSlides.Range("F18").Value = "P2"
$wsSlides.Range("H18").Value = "Título sem destaque em dourado; tipografia alinhada ao padrão."
$wsSlides.Range("J18").Value = "-"

# Row 21 (S51 - Metas por categoria de risco)
$wsSlides.Range("E21").Value = "Novo"
$wsSlides.Range("H21").Value = "Tabela-resumo com alvos e notas práticas"
$wsSlides.Range("J21").Value = "-"

# Row 23 (S53 - Base de evidência CTT + trials)
$wsSlides.Range("E23").Value = "Novo"
$wsSlides.Range("F23").Value = "P0"
$wsSlides.Range("H23").Value = "Resumo fundação + trials + 'living evidence'"
$wsSlides.Range("J23").Value = "-"

# Row 25 (S55 - Viés de publicação)
$wsSlides.Range("E25").Value = "Novo"
$wsSlides.Range("F25").Value = "P0"
$wsSlides.Range("H25").Value = "Sinais e ações práticas no GRADE"
$wsSlides.Range("J25").Value = "-"

# Row 26 (S56 - Divergências entre diretrizes)
$wsSlides.Range("E26").Value = "Novo"
$wsSlides.Range("H26").Value = "Tabela comparativa SBC vs ESC/EAS vs ACC vs AACE"
$wsSlides.Range("J26").Value = "-"

# Rows 30 and 31 (S08, S22) are dropped entirely; deleting row 30 twice
# removes both (the old row 31 shifts into row 30's place).
$wsSlides.Rows.Item(30).Delete()
$wsSlides.Rows.Item(30).Delete()

# --- Batches sheet ---

# Row 5 (Patch 2.7 entry) is dropped entirely.
$wsBatches.Rows.Item(5).Delete()
